$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.600.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5275"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3170"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06801"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.839.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.016"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007935"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.622.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.073.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.616"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.341"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.217"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.708"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.220"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08704"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.083"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04860"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7320"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.866"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.093"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.344"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01733"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4822"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9034"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.912"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.711"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4201"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.106"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1246"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05823"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8940"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.97%  "
